$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "correcion base de datos": shorten / correct the long "Causas" labels ---
# Row 5 ("Edema, proteinuria y trastornos hipertensivos en el embarazo, parto y el
#        puerperio") is shortened.
$ws.Range("B5").Value = "Edema, proteinuria y trastornos hipertensivos"

# Row 7 ("Complicaciones predominantes, relacionadas con el embarazo y el parto")
#        is shortened.
$ws.Range("B7").Value = "Complicaciones predominantes"

# Row 11 ("Muerte materna debida a cualquier causa obstétrica que ocurre después
#        de 42 días pero antes de un año del parto") is shortened.
$ws.Range("B11").Value = "Cualquier causa obstétrica que ocurre después de 42 días pero antes de un año del parto"

# --- Widen column B to fit the longer labels; keep C and D at their width ---
$ws.Columns.Item(2).ColumnWidth = 52.67
$ws.Columns.Item(3).ColumnWidth = 11.17
$ws.Columns.Item(4).ColumnWidth = 11.17

# --- Move / restore the active selection ---
$ws.Range("G7").Select() | Out-Null

# --- Turn on iterative calculation for the workbook ---
$excel.Iteration = $true
$excel.MaxIterations = 100
